$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row corresponding to the "「等々」" entry (row 833),
# shifting all subsequent rows up by one.
$ws.Rows.Item(833).Delete()
